$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new log rows for the latest success-rate runs (p2p and first mile)
$ws.Range("A3").Value = "2025-08-11 12:43:13"
$ws.Range("B3").Value = "CMM0246LCL0039"
$ws.Range("C3").Value = "Success"
$ws.Range("D3").Value = "not selected"
$ws.Range("G3").Value = 3

$ws.Range("A4").Value = "2025-08-11 12:46:14"
$ws.Range("B4").Value = "CMM0246LCL0039"
$ws.Range("C4").Value = "Success"
$ws.Range("D4").Value = "not selected"
$ws.Range("G4").Value = 3

$ws.Range("A5").Value = "2025-08-11 12:53:50"
$ws.Range("B5").Value = "CMM0246LCL0039"
$ws.Range("C5").Value = "Success"
$ws.Range("D5").Value = "not selected"
$ws.Range("G5").Value = 3
